# Update the "Pais" worksheet with refreshed COVID-19 country statistics
# and bump the "Datos actualizados..." timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Octubre de 2020 a las 18:54"

# Column layout: A=Pais, B=Casos totales, C=Nuevos casos, D=Casos activos,
# E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes

# Row 4 - Estados Unidos
$ws.Cells.Item(4, 2).Value = 8108686
$ws.Cells.Item(4, 3).Value = 18433
$ws.Cells.Item(4, 4).Value = 5239326
$ws.Cells.Item(4, 5).Value = 2648180
$ws.Cells.Item(4, 7).Value = 307
$ws.Cells.Item(4, 8).Value = 221180

# Row 5 - India
$ws.Cells.Item(5, 2).Value = 7301804
$ws.Cells.Item(5, 3).Value = 64722
$ws.Cells.Item(5, 4).Value = 6376863
$ws.Cells.Item(5, 5).Value = 813669
$ws.Cells.Item(5, 7).Value = 655
$ws.Cells.Item(5, 8).Value = 111272

# Row 6 - Brasil
$ws.Cells.Item(6, 2).Value = 5117825
$ws.Cells.Item(6, 3).Value = 3002
$ws.Cells.Item(6, 5).Value = 439689
$ws.Cells.Item(6, 7).Value = 98
$ws.Cells.Item(6, 8).Value = 151161

# Row 8 - España
$ws.Cells.Item(8, 2).Value = 937311
$ws.Cells.Item(8, 3).Value = 11970
$ws.Cells.Item(8, 7).Value = 209
$ws.Cells.Item(8, 8).Value = 33413

# Row 27 - Israel
$ws.Cells.Item(27, 2).Value = 298500
$ws.Cells.Item(27, 3).Value = 1848
$ws.Cells.Item(27, 4).Value = 251711
$ws.Cells.Item(27, 5).Value = 44691
$ws.Cells.Item(27, 7).Value = 43
$ws.Cells.Item(27, 8).Value = 2098

# Row 30 - Canada
$ws.Cells.Item(30, 2).Value = 188805
$ws.Cells.Item(30, 3).Value = 1924
$ws.Cells.Item(30, 4).Value = 159018
$ws.Cells.Item(30, 5).Value = 20126
$ws.Cells.Item(30, 7).Value = 7
$ws.Cells.Item(30, 8).Value = 9661

# Row 37 - Chequia
$ws.Cells.Item(37, 2).Value = 135425
$ws.Cells.Item(37, 3).Value = 5678
$ws.Cells.Item(37, 4).Value = 60804
$ws.Cells.Item(37, 5).Value = 73463
$ws.Cells.Item(37, 7).Value = 52
$ws.Cells.Item(37, 8).Value = 1158

# Row 40 - Republica Dominicana
$ws.Cells.Item(40, 2).Value = 119662
$ws.Cells.Item(40, 3).Value = 654
$ws.Cells.Item(40, 4).Value = 95460
$ws.Cells.Item(40, 5).Value = 22016
$ws.Cells.Item(40, 7).Value = 3
$ws.Cells.Item(40, 8).Value = 2186

# Row 58 - Suiza
$ws.Cells.Item(58, 5).Value = 16795
$ws.Cells.Item(58, 7).Value = 4
$ws.Cells.Item(58, 8).Value = 2109

# Row 74 - Kenia
$ws.Cells.Item(74, 2).Value = 42541
$ws.Cells.Item(74, 3).Value = 604
$ws.Cells.Item(74, 4).Value = 31428
$ws.Cells.Item(74, 5).Value = 10316
$ws.Cells.Item(74, 7).Value = 10
$ws.Cells.Item(74, 8).Value = 797

# Row 87 - Grecia
$ws.Cells.Item(87, 2).Value = 23495
$ws.Cells.Item(87, 3).Value = 435
$ws.Cells.Item(87, 5).Value = 13037
$ws.Cells.Item(87, 7).Value = 7
$ws.Cells.Item(87, 8).Value = 469

# Row 113 - Haiti
$ws.Cells.Item(113, 2).Value = 8908
$ws.Cells.Item(113, 3).Value = 21
$ws.Cells.Item(113, 4).Value = 7182
$ws.Cells.Item(113, 5).Value = 1495
$ws.Cells.Item(113, 7).Value = 1
$ws.Cells.Item(113, 8).Value = 231

# Row 160 - Republica de Chipre
$ws.Cells.Item(160, 2).Value = 2181
$ws.Cells.Item(160, 3).Value = 51
$ws.Cells.Item(160, 5).Value = 712

# Row 185 - Isla de Man
$ws.Cells.Item(185, 2).Value = 348
$ws.Cells.Item(185, 3).Value = 2
$ws.Cells.Item(185, 4).Value = 319
